$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns that use the numeric style (s="1", numFmtId 4) seen on existing data rows.
$s1cols = @("AH","AI","AQ","AW","AX","AY","BA","BC","BD","BE","BI","BJ","BL","BN")
$templateRow = 46
$newRows = 47..52

# Copy the numeric-style formatting (only, via PasteSpecial formats) from the
# last existing data row onto each new row, cell by cell, so only the same
# columns that carry formatting on row 46 get a style on the new rows.
foreach ($r in $newRows) {
    foreach ($col in $s1cols) {
        $ws.Range("${col}${templateRow}").Copy()
        $ws.Range("${col}${r}").PasteSpecial(-4122)
    }
}
$excel.CutCopyMode = $false

# Row 47: LUIZ FERNANDO DE OLIVEITA CAETANO
$ws.Range("A47").Value = 49
$ws.Range("B47").Value = "LUIZ FERNANDO DE OLIVEITA CAETANO"
$ws.Range("C47").Value = "FERNANDINHO"
$ws.Range("AH47").Value = 0
$ws.Range("AI47").Value = 0
$ws.Range("AQ47").Value = 0.08
$ws.Range("AU47").Value = 220
$ws.Range("AV47").Value = 220
$ws.Range("AW47").Value = 0
$ws.Range("AX47").Value = 0
$ws.Range("AY47").Value = 0
$ws.Range("AZ47").Value = 0
$ws.Range("BA47").Value = 0
$ws.Range("BB47").Value = 0
$ws.Range("BC47").Value = 0
$ws.Range("BD47").Value = 0
$ws.Range("BE47").Value = 0
$ws.Range("BH47").Value = $false
$ws.Range("BI47").Value = 0
$ws.Range("BJ47").Value = 0
$ws.Range("BK47").Value = 0
$ws.Range("BL47").Value = 0
$ws.Range("BN47").Value = 0
$ws.Range("BP47").Value = 0
$ws.Range("BQ47").Value = $false

# Row 48: JANILSON DOS SANTOS
$ws.Range("A48").Value = 50
$ws.Range("B48").Value = "JANILSON DOS SANTOS"
$ws.Range("C48").Value = "GEVÃO"
$ws.Range("AH48").Value = 0
$ws.Range("AI48").Value = 0
$ws.Range("AQ48").Value = 0.08
$ws.Range("AU48").Value = 220
$ws.Range("AV48").Value = 220
$ws.Range("AW48").Value = 0
$ws.Range("AX48").Value = 0
$ws.Range("AY48").Value = 0
$ws.Range("AZ48").Value = 0
$ws.Range("BA48").Value = 0
$ws.Range("BB48").Value = 0
$ws.Range("BC48").Value = 0
$ws.Range("BD48").Value = 0
$ws.Range("BE48").Value = 0
$ws.Range("BH48").Value = $false
$ws.Range("BI48").Value = 0
$ws.Range("BJ48").Value = 0
$ws.Range("BK48").Value = 0
$ws.Range("BL48").Value = 0
$ws.Range("BN48").Value = 0
$ws.Range("BP48").Value = 0
$ws.Range("BQ48").Value = $false

# Row 49: GEOVAN DOS SANTOS
$ws.Range("A49").Value = 51
$ws.Range("B49").Value = "GEOVAN DOS SANTOS"
$ws.Range("C49").Value = "GEOVAN"
$ws.Range("AH49").Value = 0
$ws.Range("AI49").Value = 0
$ws.Range("AQ49").Value = 0.08
$ws.Range("AU49").Value = 220
$ws.Range("AV49").Value = 220
$ws.Range("AW49").Value = 0
$ws.Range("AX49").Value = 0
$ws.Range("AY49").Value = 0
$ws.Range("AZ49").Value = 0
$ws.Range("BA49").Value = 0
$ws.Range("BB49").Value = 0
$ws.Range("BC49").Value = 0
$ws.Range("BD49").Value = 0
$ws.Range("BE49").Value = 0
$ws.Range("BH49").Value = $false
$ws.Range("BI49").Value = 0
$ws.Range("BJ49").Value = 0
$ws.Range("BK49").Value = 0
$ws.Range("BL49").Value = 0
$ws.Range("BN49").Value = 0
$ws.Range("BP49").Value = 0
$ws.Range("BQ49").Value = $false

# Row 50: ERCI CARLOS PEREIRA
$ws.Range("A50").Value = 52
$ws.Range("B50").Value = "ERCI CARLOS PEREIRA"
$ws.Range("C50").Value = "ERCI"
$ws.Range("AH50").Value = 0
$ws.Range("AI50").Value = 0
$ws.Range("AQ50").Value = 0.08
$ws.Range("AU50").Value = 220
$ws.Range("AV50").Value = 220
$ws.Range("AW50").Value = 0
$ws.Range("AX50").Value = 0
$ws.Range("AY50").Value = 0
$ws.Range("AZ50").Value = 0
$ws.Range("BA50").Value = 0
$ws.Range("BB50").Value = 0
$ws.Range("BC50").Value = 0
$ws.Range("BD50").Value = 0
$ws.Range("BE50").Value = 0
$ws.Range("BH50").Value = $false
$ws.Range("BI50").Value = 0
$ws.Range("BJ50").Value = 0
$ws.Range("BK50").Value = 0
$ws.Range("BL50").Value = 0
$ws.Range("BN50").Value = 0
$ws.Range("BP50").Value = 0
$ws.Range("BQ50").Value = $false

# Row 51: EDUARDO DE ALMEIDA DOS SANTOS
$ws.Range("A51").Value = 53
$ws.Range("B51").Value = "EDUARDO DE ALMEIDA DOS SANTOS"
$ws.Range("C51").Value = "EDU"
$ws.Range("AH51").Value = 0
$ws.Range("AI51").Value = 0
$ws.Range("AQ51").Value = 0.08
$ws.Range("AU51").Value = 220
$ws.Range("AV51").Value = 220
$ws.Range("AW51").Value = 0
$ws.Range("AX51").Value = 0
$ws.Range("AY51").Value = 0
$ws.Range("AZ51").Value = 0
$ws.Range("BA51").Value = 0
$ws.Range("BB51").Value = 0
$ws.Range("BC51").Value = 0
$ws.Range("BD51").Value = 0
$ws.Range("BE51").Value = 0
$ws.Range("BH51").Value = $false
$ws.Range("BI51").Value = 0
$ws.Range("BJ51").Value = 0
$ws.Range("BK51").Value = 0
$ws.Range("BL51").Value = 0
$ws.Range("BN51").Value = 0
$ws.Range("BP51").Value = 0
$ws.Range("BQ51").Value = $false

# Row 52: ANA CAROLINA DE ALBUQUERQUE PRADO
$ws.Range("A52").Value = 54
$ws.Range("B52").Value = "ANA CAROLINA DE ALBUQUERQUE PRADO"
$ws.Range("C52").Value = "CAROL"
$ws.Range("AH52").Value = 0
$ws.Range("AI52").Value = 0
$ws.Range("AQ52").Value = 0.08
$ws.Range("AU52").Value = 220
$ws.Range("AV52").Value = 220
$ws.Range("AW52").Value = 0
$ws.Range("AX52").Value = 0
$ws.Range("AY52").Value = 0
$ws.Range("AZ52").Value = 0
$ws.Range("BA52").Value = 0
$ws.Range("BB52").Value = 0
$ws.Range("BC52").Value = 0
$ws.Range("BD52").Value = 0
$ws.Range("BE52").Value = 0
$ws.Range("BH52").Value = $false
$ws.Range("BI52").Value = 0
$ws.Range("BJ52").Value = 0
$ws.Range("BK52").Value = 0
$ws.Range("BL52").Value = 0
$ws.Range("BN52").Value = 0
$ws.Range("BP52").Value = 0
$ws.Range("BQ52").Value = $false

# Update the defined name / used range to cover the new rows (A1:BQ52).
$wb.Names("Funcionários").RefersTo = "='Funcionários'!`$A`$1:`$BQ`$52"

Write-Output "done"
